# chore: erase semester_aktif on create and export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Aktif (true/false)" header text from D1 (keep its style)
$ws.Range("D1").Value = $null

# Remove the "true" value from D2 entirely
$ws.Range("D2").Value = $null

# Reset column D back to a near-default width, now that it no longer
# needs to fit "Aktif (true/false)" / "true"
$ws.Columns("D").ColumnWidth = 8.3
